# Trade #24 closed at 2026-02-17 12:37:28 - unknown UNKNOWN +0.000%
#
# This script applies the results of closing trade #24 to the workbook:
#  - Summary sheet: updated aggregate stats
#  - Strategy Status sheet: updated MarketMaking strategy row
#  - All Trades / MarketMaking sheets: append the new trade row (#24 -> row 25)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.69               # Current Capital
$summary.Range("B4").Value = 0.6899999999999999    # Total P&L $
$summary.Range("B5").Value = 0.57                  # Total P&L %
$summary.Range("B6").Value = 24                    # Total Trades
$summary.Range("B7").Value = 10                    # Winning Trades
$summary.Range("B9").Value = 41.67                 # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.69                 # Capital
$status.Range("D4").Value = 24                     # Trades
$status.Range("E4").Value = 0.6899999999999999     # P&L $
$status.Range("F4").Value = 0.6899999999999999     # P&L %
$status.Range("G4").Value = 41.67                  # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append new trade row (#24) to both "All Trades" and "MarketMaking"
#    sheets at row 25.
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($name in $tradeSheets) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A25").Value = 24

    # Date/time columns must stay plain text, not get auto-converted to
    # Excel date/time serial numbers.
    $ws.Range("B25").NumberFormat = "@"
    $ws.Range("B25").Value = "2026-02-17"
    $ws.Range("B25").Style = "Normal"

    $ws.Range("C25").Value = "12:37:21"

    $ws.Range("D25").Value = "MarketMaking"
    $ws.Range("E25").Value = "DOWN"
    $ws.Range("F25").Value = 0.34
    $ws.Range("G25").Value = 0.39
    $ws.Range("H25").Value = "CLOSED"
    $ws.Range("I25").Value = 14.7059
    $ws.Range("J25").Value = 0.05
    $ws.Range("K25").Value = 100.69
    $ws.Range("L25").Value = 0
    $ws.Range("M25").Value = 0
    $ws.Range("N25").Value = 0.6
    $ws.Range("O25").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P25").Value = "early_exit"
    $ws.Range("Q25").Value = 0.13
}
